# #5: property aircraft done
# The "建物" (building) sheet had its "category" column (I) incorrectly
# populated with the literal "land" for every data row (rows 2-10). This
# fixes it so the building sheet reports "building" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 9).Value = "building"
}
